# "Use 'Virus Strain' column"
#
# The sheet had a "Virus Strain" column (AA) populated with generic strain
# names, plus two extra helper columns "Taxon Virus Strain" (AD) and
# "Taxon ID" (AE) that held the actual per-row taxon-qualified virus
# strain string (and a taxon id number) used for NCBI taxon lookups.
#
# This edit folds the AD "Taxon Virus Strain" values into the AA "Virus
# Strain" column (the column that's actually used/consumed), then drops
# the now-redundant AD:AE helper columns entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 13

for ($r = 4; $r -le $lastRow; $r++) {
    $taxonStrain = $ws.Cells.Item($r, 30).Value   # column AD = "Taxon Virus Strain"
    $ws.Cells.Item($r, 27).Value = $taxonStrain    # column AA = "Virus Strain"
}

# Drop the helper columns AD (Taxon Virus Strain) and AE (Taxon ID) - in
# that order so the AE delete doesn't need re-aiming after AD shifts left.
$ws.Range("AE1:AE1048576").EntireColumn.Delete()
$ws.Range("AD1:AD1048576").EntireColumn.Delete()

# Widen the Virus Strain column now that it holds the longer taxon strings.
$ws.Range("AA1").EntireColumn.ColumnWidth = 25.6640625

# Restore the row heights for the data rows (no longer need the extra
# pixel that the taxon columns' thick border used to force).
for ($r = 4; $r -le $lastRow; $r++) {
    $ws.Rows.Item($r).RowHeight = 24
}

# Reset the view: scroll/selection no longer needs to park on the removed
# taxon columns.
$ws.Application.ActiveWindow.ScrollColumn = 21
$ws.Range("AD1:AE1048576").Select()
